# Update cryptos list - applies latest scraped values to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $cell = $ws.Range($CellRef)
    # Force any numeric-looking strings (e.g. "213.62") to be stored as
    # literal text, matching the original inline-string cell type, then
    # strip the temporary number-format style so no stray style survives.
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "27.939.46"
Set-TextValue "E2" "  +1.49%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.642.54"
Set-TextValue "E3" "  +1.20%  "

# Row 5 - BNB
Set-TextValue "D5" "213.62"
Set-TextValue "E5" "  +0.90%  "

# Row 6 - XRP
Set-TextValue "E6" "  -0.08%  "

# Row 7
Set-TextValue "E7" "  -0.04%  "

# Row 8
Set-TextValue "D8" "23.95"
Set-TextValue "E8" "  +3.52%  "

# Row 9
Set-TextValue "E9" "  +0.66%  "

# Row 10
Set-TextValue "E10" "  +0.84%  "

# Row 11
Set-TextValue "D11" "0.0876"
Set-TextValue "E11" "  -0.42%  "

# Row 12
Set-TextValue "D12" "1.875.74"
Set-TextValue "E12" "  +1.20%  "

# Row 13
Set-TextValue "D13" "1.646.28"
Set-TextValue "E13" "  +1.47%  "

# Row 14
Set-TextValue "E14" "  +4.98%  "

# Row 15
Set-TextValue "E15" "  +1.08%  "

# Row 16
Set-TextValue "D16" "65.97"
Set-TextValue "E16" "  +1.23%  "

# Row 17
Set-TextValue "D17" "27.917.36"
Set-TextValue "E17" "  +1.55%  "

# Row 18
Set-TextValue "D18" "230.97"
Set-TextValue "E18" "  +0.47%  "

# Row 19 - ShibaInu (subscript three character U+2083)
Set-TextValue "D19" "0.0$([char]0x2083)0727"
Set-TextValue "E19" "  +1.15%  "

# Row 20
Set-TextValue "E20" "  +1.44%  "

# Row 21
Set-TextValue "E21" "  -0.06%  "

# Row 22
Set-TextValue "D22" "11.12"
Set-TextValue "E22" "  +7.42%  "

# Row 23
Set-TextValue "E23" "  +1.66%  "

# Row 24
Set-TextValue "E24" "  -0.19%  "

# Row 25
Set-TextValue "D25" "152.34"
Set-TextValue "E25" "  +2.36%  "

# Row 26
Set-TextValue "E26" "  +1.00%  "

# Row 27
Set-TextValue "E27" "  +0.91%  "

# Row 28
Set-TextValue "E28" "  +1.40%  "

# Row 29
Set-TextValue "E29" "  -0.02%  "

# Row 31
Set-TextValue "E31" "  +0.51%  "

# Row 32
Set-TextValue "E32" "  +2.22%  "

# Row 33 - swaps with row 34: becomes InternetComputer(DFINITY)
Set-TextValue "B33" "InternetComputer(DFINITY)"
Set-TextValue "C33" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D33" "3.12"
Set-TextValue "E33" "  +2.32%  "

# Row 34 - becomes Maker
Set-TextValue "B34" "Maker"
Set-TextValue "C34" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D34" "1.423.03"
Set-TextValue "E34" "  -3.16%  "

# Row 35
Set-TextValue "E35" "  +2.28%  "

# Row 36
Set-TextValue "E36" "  +0.24%  "

# Row 37
Set-TextValue "D37" "0.893"
Set-TextValue "E37" "  +2.47%  "

# Row 38
Set-TextValue "E38" "  -0.77%  "

# Row 39
Set-TextValue "E39" "  +1.06%  "

# Row 40
Set-TextValue "E40" "  +0.81%  "

# Row 41
Set-TextValue "E41" "  +2.13%  "

# Row 43
Set-TextValue "D43" "67.26"
Set-TextValue "E43" "  +0.12%  "

# Row 44
Set-TextValue "E44" "  +0.50%  "

# Row 45
Set-TextValue "D45" "5.45"
Set-TextValue "E45" "  +3.18%  "

# Row 46
Set-TextValue "E46" "  +3.65%  "

# Row 47
Set-TextValue "E47" "  +0.37%  "

# Row 48
Set-TextValue "D48" "1.784.26"
Set-TextValue "E48" "  +1.20%  "

# Row 49
Set-TextValue "D49" "88.89"
Set-TextValue "E49" "  +1.81%  "

# Row 50
Set-TextValue "E50" "  +1.05%  "

# Row 51
Set-TextValue "D51" "0.0507"
Set-TextValue "E51" "  +0.67%  "
